$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue ($ws.Range("D2")) '46.470.66'
Set-TextValue ($ws.Range("E2")) '  -0.39%  '

# Row 3
Set-TextValue ($ws.Range("D3")) '2.424.95'
Set-TextValue ($ws.Range("E3")) '  +6.64%  '

# Row 4
Set-TextValue ($ws.Range("D4")) '1.00'
Set-TextValue ($ws.Range("E4")) '  +0.03%  '

# Row 5
Set-TextValue ($ws.Range("D5")) '296.65'
Set-TextValue ($ws.Range("E5")) '  -1.64%  '

# Row 6
Set-TextValue ($ws.Range("D6")) '97.97'
Set-TextValue ($ws.Range("E6")) '  -2.43%  '

# Row 7
Set-TextValue ($ws.Range("E7")) '  +0.85%  '

# Row 8
Set-TextValue ($ws.Range("D8")) '1.00'
Set-TextValue ($ws.Range("E8")) '  +0.23%  '

# Row 9
Set-TextValue ($ws.Range("D9")) '0.512'
Set-TextValue ($ws.Range("E9")) '  +0.60%  '

# Row 10
Set-TextValue ($ws.Range("D10")) '35.41'
Set-TextValue ($ws.Range("E10")) '  +0.48%  '

# Row 11
Set-TextValue ($ws.Range("D11")) '0.0788'
Set-TextValue ($ws.Range("E11")) '  -1.68%  '

# Row 12
Set-TextValue ($ws.Range("D12")) '7.15'
Set-TextValue ($ws.Range("E12")) '  +0.89%  '

# Row 13
Set-TextValue ($ws.Range("D13")) '0.105'
Set-TextValue ($ws.Range("E13")) '  +2.14%  '

# Row 14
Set-TextValue ($ws.Range("D14")) '2.797.78'
Set-TextValue ($ws.Range("E14")) '  +6.81%  '

# Row 15
Set-TextValue ($ws.Range("D15")) '2.434.26'
Set-TextValue ($ws.Range("E15")) '  +7.05%  '

# Row 16
Set-TextValue ($ws.Range("D16")) '0.851'
Set-TextValue ($ws.Range("E16")) '  +6.30%  '

# Row 17
Set-TextValue ($ws.Range("D17")) '14.08'
Set-TextValue ($ws.Range("E17")) '  +2.91%  '

# Row 18
Set-TextValue ($ws.Range("D18")) '46.416.80'
Set-TextValue ($ws.Range("E18")) '  -0.48%  '

# Row 19
Set-TextValue ($ws.Range("D19")) '12.75'
Set-TextValue ($ws.Range("E19")) '  +0.52%  '

# Row 20
Set-TextValue ($ws.Range("D20")) '0.0₃0951'
Set-TextValue ($ws.Range("E20")) '  -1.86%  '

# Row 21
Set-TextValue ($ws.Range("D21")) '6.24'
Set-TextValue ($ws.Range("E21")) '  +6.48%  '

# Row 22
Set-TextValue ($ws.Range("D22")) '67.50'
Set-TextValue ($ws.Range("E22")) '  +2.25%  '

# Row 23
Set-TextValue ($ws.Range("D23")) '245.35'
Set-TextValue ($ws.Range("E23")) '  -1.57%  '

# Row 24
Set-TextValue ($ws.Range("E24")) '  +0.43%  '

# Row 25
Set-TextValue ($ws.Range("D25")) '1.96'
Set-TextValue ($ws.Range("E25")) '  +4.76%  '

# Row 26
Set-TextValue ($ws.Range("E26")) '  -0.08%  '

# Row 27
Set-TextValue ($ws.Range("D27")) '39.57'
Set-TextValue ($ws.Range("E27")) '  -4.76%  '

# Row 28
Set-TextValue ($ws.Range("D28")) '2.22'
Set-TextValue ($ws.Range("E28")) '  -1.69%  '

# Row 29
Set-TextValue ($ws.Range("D29")) '9.80'
Set-TextValue ($ws.Range("E29")) '  +1.56%  '

# Row 30
Set-TextValue ($ws.Range("E30")) '  +13.93%  '

# Row 31
Set-TextValue ($ws.Range("D31")) '21.39'
Set-TextValue ($ws.Range("E31")) '  +5.78%  '

# Row 32
Set-TextValue ($ws.Range("E32")) '  -1.75%  '

# Row 33
Set-TextValue ($ws.Range("D33")) '5.56'
Set-TextValue ($ws.Range("E33")) '  +3.08%  '

# Row 34
Set-TextValue ($ws.Range("D34")) '148.05'
Set-TextValue ($ws.Range("E34")) '  +0.61%  '

# Row 35
Set-TextValue ($ws.Range("D35")) '0.0772'
Set-TextValue ($ws.Range("E35")) '  -0.16%  '

# Row 36
Set-TextValue ($ws.Range("D36")) '1.99'
Set-TextValue ($ws.Range("E36")) '  +17.09%  '

# Row 37
Set-TextValue ($ws.Range("E37")) '  +0.72%  '

# Row 38
Set-TextValue ($ws.Range("E38")) '  +0.23%  '

# Row 39
Set-TextValue ($ws.Range("D39")) '15.33'
Set-TextValue ($ws.Range("E39")) '  -3.13%  '

# Row 40
Set-TextValue ($ws.Range("D40")) '3.97'
Set-TextValue ($ws.Range("E40")) '  +2.24%  '

# Row 41
Set-TextValue ($ws.Range("D41")) '0.0303'
Set-TextValue ($ws.Range("E41")) '  +2.26%  '

# Row 42
Set-TextValue ($ws.Range("D42")) '3.26'
Set-TextValue ($ws.Range("E42")) '  +4.08%  '

# Row 43
Set-TextValue ($ws.Range("D43")) '1.978.86'
Set-TextValue ($ws.Range("E43")) '  +10.33%  '

# Row 44
Set-TextValue ($ws.Range("D44")) '1.00'
Set-TextValue ($ws.Range("E44")) '  +0.19%  '

# Row 45
Set-TextValue ($ws.Range("D45")) '92.23'
Set-TextValue ($ws.Range("E45")) '  -1.05%  '

# Row 46
Set-TextValue ($ws.Range("E46")) '  -2.18%  '

# Row 47
Set-TextValue ($ws.Range("D47")) '16.30'
Set-TextValue ($ws.Range("E47")) '  +30.24%  '

# Row 48
Set-TextValue ($ws.Range("D48")) '8.64'
Set-TextValue ($ws.Range("E48")) '  +9.19%  '

# Row 49
Set-TextValue ($ws.Range("D49")) '101.07'
Set-TextValue ($ws.Range("E49")) '  +6.28%  '

# Row 50
Set-TextValue ($ws.Range("D50")) '2.670.46'
Set-TextValue ($ws.Range("E50")) '  +7.01%  '

# Row 51
Set-TextValue ($ws.Range("D51")) '0.186'
Set-TextValue ($ws.Range("E51")) '  -0.13%  '
